$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Harvard case classification: new "average_doctor" column inserted at BP,
# pushing the previous average_doctor data into BQ (relabelled "average_doctor_old").

$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"
$ws.Range("E4").Value = 0.475
$ws.Range("F4").Value = 0.049
$ws.Range("G4").Value = 0.221
$ws.Range("N4").Value = 0.48
$ws.Range("O4").Value = 0.067
$ws.Range("P4").Value = 0.258
$ws.Range("Q4").Value = 0.039
$ws.Range("R4").Value = 0.028
$ws.Range("S4").Value = 0.167
$ws.Range("W4").Value = 0.365
$ws.Range("AI4").Value = 0.442
$ws.Range("AJ4").Value = 0.08500000000000001
$ws.Range("AK4").Value = 0.291
$ws.Range("AU4").Value = 0.243
$ws.Range("BA4").Value = 2.089
$ws.Range("BB4").Value = 0.142
$ws.Range("BC4").Value = 0.377
$ws.Range("BG4").Value = 0.739
$ws.Range("BH4").Value = 0.136
$ws.Range("BI4").Value = 0.369
$ws.Range("BM4").Value = 0.756
$ws.Range("BN4").Value = 0.064
$ws.Range("BO4").Value = 0.253
$ws.Range("BP4").Value = 0.696
$ws.Range("BQ4").Value = 0.765
$ws.Range("E5").Value = 0.601
$ws.Range("F5").Value = 0.051
$ws.Range("G5").Value = 0.226
$ws.Range("N5").Value = 0.703
$ws.Range("O5").Value = 0.079
$ws.Range("P5").Value = 0.281
$ws.Range("Q5").Value = 0.02
$ws.Range("R5").Value = 0.005
$ws.Range("S5").Value = 0.07099999999999999
$ws.Range("W5").Value = 0.322
$ws.Range("X5").Value = 0.092
$ws.Range("Y5").Value = 0.304
$ws.Range("AI5").Value = 0.443
$ws.Range("AJ5").Value = 0.081
$ws.Range("AK5").Value = 0.284
$ws.Range("AU5").Value = 0.448
$ws.Range("AV5").Value = 0.078
$ws.Range("AW5").Value = 0.28
$ws.Range("BA5").Value = 1.295
$ws.Range("BB5").Value = 0.066
$ws.Range("BC5").Value = 0.257
$ws.Range("BG5").Value = 0.382
$ws.Range("BH5").Value = 0.049
$ws.Range("BI5").Value = 0.22
$ws.Range("BM5").Value = 0.519
$ws.Range("BN5").Value = 0.044
$ws.Range("BO5").Value = 0.209
$ws.Range("BP5").Value = 0.432
$ws.Range("BQ5").Value = 0.455
$ws.Range("E6").Value = 0.531
$ws.Range("N6").Value = 0.57
$ws.Range("Q6").Value = 0.026
$ws.Range("W6").Value = 0.342
$ws.Range("AI6").Value = 0.442
$ws.Range("AU6").Value = 0.315
$ws.Range("BA6").Value = 1.593
$ws.Range("BG6").Value = 0.504
$ws.Range("BM6").Value = 0.615
$ws.Range("BP6").Value = 0.531
$ws.Range("BQ6").Value = 0.5679999999999999
$ws.Range("E7").Value = 0.571
$ws.Range("N7").Value = 0.643
$ws.Range("Q7").Value = 0.022
$ws.Range("W7").Value = 0.33
$ws.Range("AI7").Value = 0.443
$ws.Range("AU7").Value = 0.383
$ws.Range("BA7").Value = 1.399
$ws.Range("BG7").Value = 0.423
$ws.Range("BM7").Value = 0.554
$ws.Range("BP7").Value = 0.466
$ws.Range("BQ7").Value = 0.494
$ws.Range("E8").Value = 0.709
$ws.Range("F8").Value = 0.061
$ws.Range("G8").Value = 0.248
$ws.Range("N8").Value = 0.791
$ws.Range("O8").Value = 0.062
$ws.Range("P8").Value = 0.25
$ws.Range("Q8").Value = 0.021
$ws.Range("S8").Value = 0.108
$ws.Range("W8").Value = 0.4
$ws.Range("X8").Value = 0.12
$ws.Range("Y8").Value = 0.346
$ws.Range("AI8").Value = 0.52
$ws.Range("AJ8").Value = 0.126
$ws.Range("AK8").Value = 0.354
$ws.Range("AU8").Value = 0.403
$ws.Range("AV8").Value = 0.08599999999999999
$ws.Range("AW8").Value = 0.292
$ws.Range("BA8").Value = 1.782
$ws.Range("BB8").Value = 0.102
$ws.Range("BC8").Value = 0.319
$ws.Range("BG8").Value = 0.5610000000000001
$ws.Range("BH8").Value = 0.107
$ws.Range("BI8").Value = 0.328
$ws.Range("BM8").Value = 0.6840000000000001
$ws.Range("BN8").Value = 0.057
$ws.Range("BO8").Value = 0.238
$ws.Range("BP8").Value = 0.594
$ws.Range("BQ8").Value = 0.627
$ws.Range("E9").Value = 0.667
$ws.Range("F9").Value = 0.222
$ws.Range("G9").Value = 0.471
$ws.Range("N9").Value = 0.718
$ws.Range("O9").Value = 0.202
$ws.Range("P9").Value = 0.45
$ws.Range("W9").Value = 0.308
$ws.Range("X9").Value = 0.213
$ws.Range("Y9").Value = 0.462
$ws.Range("AI9").Value = 0.462
$ws.Range("AJ9").Value = 0.249
$ws.Range("AK9").Value = 0.499
$ws.Range("BA9").Value = 1.744
$ws.Range("BB9").Value = 0.25
$ws.Range("BC9").Value = 0.5
$ws.Range("BG9").Value = 0.59
$ws.Range("BH9").Value = 0.242
$ws.Range("BI9").Value = 0.492
$ws.Range("BM9").Value = 0.667
$ws.Range("BN9").Value = 0.222
$ws.Range("BO9").Value = 0.471
$ws.Range("BP9").Value = 0.581
$ws.Range("BQ9").Value = 0.613
$ws.Range("E10").Value = 0.821
$ws.Range("F10").Value = 0.147
$ws.Range("G10").Value = 0.384
$ws.Range("N10").Value = 0.923
$ws.Range("O10").Value = 0.07099999999999999
$ws.Range("P10").Value = 0.266
$ws.Range("W10").Value = 0.513
$ws.Range("X10").Value = 0.25
$ws.Range("Y10").Value = 0.5
$ws.Range("AI10").Value = 0.5639999999999999
$ws.Range("AJ10").Value = 0.246
$ws.Range("AK10").Value = 0.496
$ws.Range("AU10").Value = 0.41
$ws.Range("AV10").Value = 0.242
$ws.Range("AW10").Value = 0.492
$ws.Range("BA10").Value = 2.231
$ws.Range("BB10").Value = 0.202
$ws.Range("BC10").Value = 0.45
$ws.Range("BG10").Value = 0.667
$ws.Range("BH10").Value = 0.222
$ws.Range("BI10").Value = 0.471
$ws.Range("BM10").Value = 0.846
$ws.Range("BN10").Value = 0.13
$ws.Range("BO10").Value = 0.361
$ws.Range("BP10").Value = 0.744
$ws.Range("BQ10").Value = 0.773
$ws.Range("E11").Value = 0.872
$ws.Range("F11").Value = 0.112
$ws.Range("G11").Value = 0.334
$ws.Range("N11").Value = 0.923
$ws.Range("O11").Value = 0.07099999999999999
$ws.Range("P11").Value = 0.266
$ws.Range("W11").Value = 0.513
$ws.Range("X11").Value = 0.25
$ws.Range("Y11").Value = 0.5
$ws.Range("AI11").Value = 0.641
$ws.Range("AJ11").Value = 0.23
$ws.Range("AK11").Value = 0.48
$ws.Range("AU11").Value = 0.5639999999999999
$ws.Range("AV11").Value = 0.246
$ws.Range("AW11").Value = 0.496
$ws.Range("BA11").Value = 2.231
$ws.Range("BB11").Value = 0.202
$ws.Range("BC11").Value = 0.45
$ws.Range("BG11").Value = 0.667
$ws.Range("BH11").Value = 0.222
$ws.Range("BI11").Value = 0.471
$ws.Range("BM11").Value = 0.846
$ws.Range("BN11").Value = 0.13
$ws.Range("BO11").Value = 0.361
$ws.Range("BP11").Value = 0.744
$ws.Range("BQ11").Value = 0.78
$ws.Range("E12").Value = 1.441
$ws.Range("F12").Value = 0.894
$ws.Range("G12").Value = 0.945
$ws.Range("N12").Value = 1.278
$ws.Range("O12").Value = 0.312
$ws.Range("P12").Value = 0.5580000000000001
$ws.Range("W12").Value = 1.5
$ws.Range("X12").Value = 0.45
$ws.Range("Y12").Value = 0.671
$ws.Range("AI12").Value = 1.6
$ws.Range("AJ12").Value = 1.44
$ws.Range("AK12").Value = 1.2
$ws.Range("AU12").Value = 2.833
$ws.Range("AV12").Value = 3.556
$ws.Range("AW12").Value = 1.886
$ws.Range("BA12").Value = 3.825
$ws.Range("BB12").Value = 0.459
$ws.Range("BC12").Value = 0.678
$ws.Range("BG12").Value = 1.154
$ws.Range("BH12").Value = 0.207
$ws.Range("BI12").Value = 0.455
$ws.Range("BM12").Value = 1.242
$ws.Range("BN12").Value = 0.244
$ws.Range("BO12").Value = 0.494
$ws.Range("BP12").Value = 1.275
$ws.Range("BQ12").Value = 1.265
$ws.Range("E13").Value = 1.421
$ws.Range("F13").Value = 0.313
$ws.Range("G13").Value = 0.5590000000000001
$ws.Range("N13").Value = 1.709
$ws.Range("O13").Value = 0.507
$ws.Range("P13").Value = 0.712
$ws.Range("W13").Value = 0.966
$ws.Range("X13").Value = 0.192
$ws.Range("Y13").Value = 0.439
$ws.Range("AI13").Value = 1.154
$ws.Range("AJ13").Value = 0.303
$ws.Range("AK13").Value = 0.551
$ws.Range("AU13").Value = 2.024
$ws.Range("AV13").Value = 0.351
$ws.Range("AW13").Value = 0.592
$ws.Range("BA13").Value = 2.095
$ws.Range("BB13").Value = 0.262
$ws.Range("BC13").Value = 0.512
$ws.Range("BG13").Value = 0.534
$ws.Range("BH13").Value = 0.048
$ws.Range("BI13").Value = 0.22
$ws.Range("BM13").Value = 0.762
$ws.Range("BN13").Value = 0.13
$ws.Range("BO13").Value = 0.361
$ws.Range("BP13").Value = 0.698
$ws.Range("BQ13").Value = 0.651
